{"js": "// Collapse the section's per-resource/case paragraphs into the section's\n// headnote paragraph.\n//\n// The document has a \"Section Headnote\" paragraph (\"What is a corporation?\")\n// followed by a run of \"Resource Number\" / \"Resource Title\" / \"Resource\n// Headnote\" / \"Case Text\" / \"Section Number\" / \"Section Title\" paragraphs\n// that ends at the *next* \"Section Headnote\" paragraph (\"This is the second\n// chapter of the casebook.\"). All of those paragraphs (including the\n// trailing Section Headnote one) get merged away: every paragraph's text\n// EXCEPT paragraphs whose style ends in \"Headnote\" is concatenated (in\n// document order) into a single run that replaces the first paragraph's\n// text, and the rest of the paragraphs are deleted.\n\nfunction xmlEscape(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the first paragraph of the span: the \"Section Headnote\" paragraph\n// with the known starting text.\nlet startIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].style === \"Section Headnote\" && items[i].text === \"What is a corporation?\") {\n    startIndex = i;\n    break;\n  }\n}\nif (startIndex === -1) {\n  throw new Error(\"Could not find the starting Section Headnote paragraph.\");\n}\n\n// Locate the end of the span: the next \"Section Headnote\" paragraph after\n// startIndex (the closing \"This is the second chapter...\" paragraph).\nlet endIndex = -1;\nfor (let i = startIndex + 1; i < items.length; i++) {\n  if (items[i].style === \"Section Headnote\") {\n    endIndex = i;\n    break;\n  }\n}\nif (endIndex === -1) {\n  throw new Error(\"Could not find the ending Section Headnote paragraph.\");\n}\n\n// Build the merged text: every paragraph strictly between start and end\n// (inclusive of end), excluding any paragraph whose style ends with\n// \"Headnote\" (\"Resource Headnote\" / \"Section Headnote\" paragraphs\n// contribute nothing to the merged text).\nlet merged = \"\";\nfor (let i = startIndex + 1; i <= endIndex; i++) {\n  const style = items[i].style || \"\";\n  if (style.endsWith(\"Headnote\")) {\n    continue;\n  }\n  merged += items[i].text;\n}\n\n// Delete every paragraph from startIndex+1 through endIndex (inclusive),\n// leaving only the first paragraph of the span.\nfor (let i = startIndex + 1; i <= endIndex; i++) {\n  items[i].delete();\n}\nawait context.sync();\n\n// Replace the first paragraph's text with the merged text. insertText()\n// treats embedded \"\\n\" as a paragraph break, which we don't want here (the\n// target keeps everything inside one run/one <w:t>), so we use insertOoxml\n// with a minimal package instead -- that preserves a literal newline\n// character inside the text node.\nconst startRange = items[startIndex].getRange();\nconst pkg =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n  '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n  '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n  '</Relationships></pkg:xmlData></pkg:part>' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"SectionHeadnote\"/></w:pPr><w:r><w:t xml:space=\"preserve\">' +\n  xmlEscape(merged) +\n  '</w:t></w:r></w:p>' +\n  '<w:sectPr><w:pgSz w:w=\"12240\" w:h=\"15840\"/></w:sectPr>' +\n  '</w:body></w:document></pkg:xmlData></pkg:part>' +\n  '</pkg:package>';\nstartRange.insertOoxml(pkg, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Collapse the section's per-resource/case paragraphs into the section's\n# headnote paragraph.\n#\n# The document has a \"Section Headnote\" paragraph (\"What is a corporation?\")\n# followed by a run of \"Resource Number\" / \"Resource Title\" / \"Resource\n# Headnote\" / \"Case Text\" / \"Section Number\" / \"Section Title\" paragraphs\n# that ends at the *next* \"Section Headnote\" paragraph (\"This is the second\n# chapter of the casebook.\"). All of those paragraphs (including the\n# trailing Section Headnote one) get merged away: every paragraph's text\n# EXCEPT paragraphs whose style ends in \"Headnote\" is concatenated (in\n# document order) into the first paragraph's run, and the rest of the\n# paragraphs are deleted.\n\n$d = $word.ActiveDocument\n\nfunction TrimMark([string]$s) {\n    # Paragraph.Range.Text includes the trailing paragraph mark (CR, 0x0D;\n    # or cell-mark BEL, 0x07, inside tables) -- strip it before comparing\n    # or concatenating paragraph text.\n    return $s.TrimEnd([char]0x0D, [char]0x07)\n}\n\n$count = $d.Paragraphs.Count\n\n# Locate the first paragraph of the span: the \"Section Headnote\" paragraph\n# with the known starting text.\n$startIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Style.NameLocal -eq \"Section Headnote\" -and (TrimMark $p.Range.Text) -eq \"What is a corporation?\") {\n        $startIndex = $i\n        break\n    }\n}\nif ($startIndex -eq -1) {\n    throw \"Could not find the starting Section Headnote paragraph.\"\n}\n\n# Locate the end of the span: the next \"Section Headnote\" paragraph after\n# startIndex (the closing \"This is the second chapter...\" paragraph).\n$endIndex = -1\nfor ($i = $startIndex + 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Style.NameLocal -eq \"Section Headnote\") {\n        $endIndex = $i\n        break\n    }\n}\nif ($endIndex -eq -1) {\n    throw \"Could not find the ending Section Headnote paragraph.\"\n}\n\n# Build the merged text: every paragraph strictly between start and end\n# (inclusive of end), excluding any paragraph whose style ends with\n# \"Headnote\" (\"Resource Headnote\" / \"Section Headnote\" paragraphs\n# contribute nothing to the merged text).\n$merged = \"\"\nfor ($i = $startIndex + 1; $i -le $endIndex; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $style = $p.Range.Style.NameLocal\n    if ($style.EndsWith(\"Headnote\")) {\n        continue\n    }\n    $merged += (TrimMark $p.Range.Text)\n}\n\n# Delete every paragraph from startIndex+1 through endIndex (inclusive) in\n# one shot by deleting the Range spanning them.\n$spanStart = $d.Paragraphs.Item($startIndex + 1).Range.Start\n$spanEnd = $d.Paragraphs.Item($endIndex).Range.End\n$span = $d.Range($spanStart, $spanEnd)\n$span.Delete()\n\n# Replace the first paragraph's text (the paragraph mark is preserved\n# automatically since we only assign to its Range, not delete it).\n$startRange = $d.Paragraphs.Item($startIndex).Range\n$startRange.Text = $merged\n"}
